# Auto-generated edit script: updates market/profit values across all 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3850.2727
$ws.Range("J17").Value = 3379
$ws.Range("L17").Value = 10137
$ws.Range("N17").Value = -10473
$ws.Range("H112").Value = 2904.524
$ws.Range("J112").Value = 3221.9443
$ws.Range("L112").Value = 9665.832900000001
$ws.Range("N112").Value = -11881.8329
$ws.Range("H138").Value = 3463.3713
$ws.Range("I138").Value = 3435.5186
$ws.Range("J138").Value = 3480.8604
$ws.Range("K138").Value = 10306.5558
$ws.Range("L138").Value = 10442.5812
$ws.Range("M138").Value = -5166.5558
$ws.Range("N138").Value = -20722.5812

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 223169.36
$ws.Range("I2").Value = 309370.34
$ws.Range("J2").Value = 1509.7142
$ws.Range("K2").Value = 309370.34
$ws.Range("L2").Value = 1509.7142
$ws.Range("M2").Value = -309257.34
$ws.Range("N2").Value = -1735.7142
$ws.Range("H32").Value = 15818.698
$ws.Range("I32").Value = 12782.0625
$ws.Range("K32").Value = 12782.0625
$ws.Range("M32").Value = -12495.0625
$ws.Range("H74").Value = 2753.9
$ws.Range("I74").Value = 3999.5
$ws.Range("K74").Value = 3999.5
$ws.Range("M74").Value = -3125.5
$ws.Range("H77").Value = 2753.9
$ws.Range("I77").Value = 3999.5
$ws.Range("K77").Value = 19997.5
$ws.Range("M77").Value = -15629.5
$ws.Range("H110").Value = 660.3333
$ws.Range("I110").Value = 660.3333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 660.3333
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1384.6667
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 223169.36
$ws.Range("I116").Value = 309370.34
$ws.Range("J116").Value = 1509.7142
$ws.Range("K116").Value = 309370.34
$ws.Range("L116").Value = 1509.7142
$ws.Range("M116").Value = -307076.34
$ws.Range("N116").Value = -6097.7142
$ws.Range("H132").Value = 1692.9452
$ws.Range("I132").Value = 1296.7297
$ws.Range("K132").Value = 3890.189100000001
$ws.Range("M132").Value = -1360.189100000001
$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 65000
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 223169.36
$ws.Range("I3").Value = 309370.34
$ws.Range("J3").Value = 1509.7142
$ws.Range("K3").Value = 309370.34
$ws.Range("L3").Value = 1509.7142
$ws.Range("M3").Value = -309256.34
$ws.Range("N3").Value = -1737.7142
$ws.Range("H94").Value = 1074.3636
$ws.Range("I94").Value = 589.875
$ws.Range("K94").Value = 589.875
$ws.Range("M94").Value = -138.875
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 76448.75
$ws.Range("J140").Value = 76448.75
$ws.Range("L140").Value = 76448.75
$ws.Range("N140").Value = -86808.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3916.6667
$ws.Range("I76").Value = 500
$ws.Range("J76").Value = 4600
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 13800
$ws.Range("M76").Value = -1117
$ws.Range("N76").Value = -14566
$ws.Range("H79").Value = 3916.6667
$ws.Range("I79").Value = 500
$ws.Range("J79").Value = 4600
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 13800
$ws.Range("M79").Value = -174
$ws.Range("N79").Value = -16452
$ws.Range("H122").Value = 1222.0667
$ws.Range("I122").Value = 662.5
$ws.Range("J122").Value = 1425.5454
$ws.Range("K122").Value = 5962.5
$ws.Range("L122").Value = 12829.9086
$ws.Range("M122").Value = -3512.5
$ws.Range("N122").Value = -17729.9086
$ws.Range("H125").Value = 8144.2856
$ws.Range("I125").Value = 2936.6667
$ws.Range("K125").Value = 8810.000100000001
$ws.Range("M125").Value = -3890.000100000001
$ws.Range("H137").Value = 3924.1667
$ws.Range("I137").Value = 2465.077
$ws.Range("J137").Value = 5648.5454
$ws.Range("K137").Value = 7395.231000000001
$ws.Range("L137").Value = 16945.6362
$ws.Range("M137").Value = -2295.231000000001
$ws.Range("N137").Value = -27145.6362

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15027.556
$ws.Range("J92").Value = 15027.556
$ws.Range("L92").Value = 15027.556
$ws.Range("N92").Value = -18771.556
$ws.Range("H93").Value = 26700
$ws.Range("J93").Value = 26700
$ws.Range("L93").Value = 26700
$ws.Range("N93").Value = -30444
$ws.Range("H132").Value = 634477.8
$ws.Range("I132").Value = 920227.9399999999
$ws.Range("K132").Value = 2760683.82
$ws.Range("M132").Value = -2758153.82
$ws.Range("H138").Value = 75500
$ws.Range("J138").Value = 75500
$ws.Range("L138").Value = 75500
$ws.Range("N138").Value = -85780
$ws.Range("H140").Value = 49999
$ws.Range("J140").Value = 49999
$ws.Range("L140").Value = 49999
$ws.Range("N140").Value = -60359

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3000
$ws.Range("I93").Value = 3000
$ws.Range("K93").Value = 3000
$ws.Range("M93").Value = -1752
$ws.Range("H96").Value = 74997
$ws.Range("J96").Value = 74997
$ws.Range("L96").Value = 74997
$ws.Range("N96").Value = -80489
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920
$ws.Range("H132").Value = 5505.579
$ws.Range("I132").Value = 4056.1538
$ws.Range("J132").Value = 8646
$ws.Range("K132").Value = 12168.4614
$ws.Range("L132").Value = 25938
$ws.Range("M132").Value = -9638.4614
$ws.Range("N132").Value = -30998
$ws.Range("H136").Value = 4399.4
$ws.Range("I136").Value = 1999
$ws.Range("K136").Value = 5997
$ws.Range("M136").Value = -3447
$ws.Range("H139").Value = 66803.5
$ws.Range("J139").Value = 66803.5
$ws.Range("L139").Value = 66803.5
$ws.Range("N139").Value = -77083.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1883.7174
$ws.Range("I132").Value = 1716.4839
$ws.Range("J132").Value = 2229.3333
$ws.Range("K132").Value = 5149.4517
$ws.Range("L132").Value = 6687.999899999999
$ws.Range("M132").Value = -2619.4517
$ws.Range("N132").Value = -11747.9999

